$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.987.78'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.685.05'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.51'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.01'
$ws.Range('E10').Value = '  +4.12%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '1.919.91'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '1.682.26'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.533'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.91'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '27.016.71'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '8.18'
$ws.Range('E18').Value = '  +5.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '236.98'
$ws.Range('E19').Value = '  +2.32%  '
$ws.Range('D20').Value = '0.0₃0735'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.23'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.11'
$ws.Range('E24').Value = '  -4.28%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.48'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.07'
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('E28').Value = '  -3.03%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0501'
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').Value = '1.498.05'
$ws.Range('E33').Value = '  +2.96%  '
$ws.Range('E34').Value = '  +0.93%  '
$ws.Range('E35').Value = '  +5.10%  '
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.588'
$ws.Range('E37').Value = '  +4.14%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.916'
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0174'
$ws.Range('E39').Value = '  +3.96%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.04'
$ws.Range('E40').Value = '  +7.35%  '
$ws.Range('E41').Value = '  -4.65%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '67.67'
$ws.Range('E43').Value = '  +3.17%  '
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').Value = '1.824.51'
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.781'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.59'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.53'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.104'
$ws.Range('E49').Value = '  +3.70%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  +8.20%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.78'
$ws.Range('E51').Value = '  +2.71%  '
